$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 567, shifting the existing data (old rows
# 567-655) down to 569-657.
$ws.Rows("567:568").Insert()

# Populate the newly inserted row 567 (new weekly observation, same
# categorical attributes as the price point that used to sit there).
$ws.Range("A567").Value = 6
$ws.Range("B567").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C567").Value = "Metropolitana"
$ws.Range("D567").Value = 44984
$ws.Range("E567").Value = 13
$ws.Range("F567").Value = 100112052
$ws.Range("G567").Value = "Albahaca"
$ws.Range("H567").Value = "Sin especificar"
$ws.Range("I567").Value = "Primera"
$ws.Range("J567").Value = 630
$ws.Range("K567").Value = 3000
$ws.Range("L567").Value = 3500
$ws.Range("M567").Value = 3127
$ws.Range("N567").Value = "$/docena de matas"
$ws.Range("O567").Value = "Región Metropolitana"
$ws.Range("P567").Value = 521
$ws.Range("Q567").Value = 6
$ws.Range("R567").Value = "Hortaliza"

# Populate the newly inserted row 568.
$ws.Range("A568").Value = 6
$ws.Range("B568").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C568").Value = "Metropolitana"
$ws.Range("D568").Value = 44984
$ws.Range("E568").Value = 13
$ws.Range("F568").Value = 100112052
$ws.Range("G568").Value = "Albahaca"
$ws.Range("H568").Value = "Sin especificar"
$ws.Range("I568").Value = "Segunda"
$ws.Range("J568").Value = 150
$ws.Range("K568").Value = 2500
$ws.Range("L568").Value = 2500
$ws.Range("M568").Value = 2500
$ws.Range("N568").Value = "$/docena de matas"
$ws.Range("O568").Value = "Región Metropolitana"
$ws.Range("P568").Value = 417
$ws.Range("Q568").Value = 6
$ws.Range("R568").Value = "Hortaliza"
